$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6: "Different hooking methods" Done? -> Skipped (Neutral style)
$ws.Range("G6").Value = "Skipped"
$ws.Range("G6").Style = "Neutral"

# Row 10: "Improve hiding sockets" Done? -> Yes (Good style)
$ws.Range("G10").Value = "Yes"
$ws.Range("G10").Style = "Gut"

# Update selection to match final cursor position
$ws.Range("H16").Select()
